# Fix typos on slide LinkedLists ProsCons

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# Shape id=11 ("При работе с большими объектами перемещение указателей проще чем копирование")
$shape1 = $s.Shapes.Item(8)
$shape1.TextFrame.TextRange.Text = "При работе с большими объектами перемещение указателей проще, чем копирование"

# Shape id=14 ("Использую дополнительную память для указателей")
$shape2 = $s.Shapes.Item(11)
$shape2.TextFrame.TextRange.Text = "Используют дополнительную память для указателей"
